$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "45.307.15"
$ws.Range("E2").Value = "  -0.88%  "

# Row 3
$ws.Range("D3").Value = "2.367.06"
$ws.Range("E3").Value = "  -1.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'330.19"
$ws.Range("E5").Value = "  +3.90%  "

# Row 6
$ws.Range("D6").Value = "'108.14"
$ws.Range("E6").Value = "  -6.01%  "

# Row 7
$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = "  -0.29%  "

# Row 8
$ws.Range("E8").Value = "  +0.13%  "

# Row 9
$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = "  -2.37%  "

# Row 10
$ws.Range("D10").Value = "'41.14"
$ws.Range("E10").Value = "  -3.73%  "

# Row 11
$ws.Range("D11").Value = "'0.0918"
$ws.Range("E11").Value = "  -1.58%  "

# Row 12
$ws.Range("D12").Value = "'8.48"
$ws.Range("E12").Value = "  -3.42%  "

# Row 13
$ws.Range("E13").Value = "  -0.38%  "

# Row 14
$ws.Range("D14").Value = "'0.980"
$ws.Range("E14").Value = "  -3.40%  "

# Row 15
$ws.Range("D15").Value = "2.728.02"
$ws.Range("E15").Value = "  -1.65%  "

# Row 16
$ws.Range("D16").Value = "'15.41"
$ws.Range("E16").Value = "  -3.67%  "

# Row 17
$ws.Range("D17").Value = "2.366.71"
$ws.Range("E17").Value = "  -1.63%  "

# Row 18
$ws.Range("D18").Value = "45.312.50"
$ws.Range("E18").Value = "  -0.92%  "

# Row 19
$ws.Range("D19").Value = "'15.23"
$ws.Range("E19").Value = "  +12.09%  "

# Row 20
$ws.Range("D20").Value = "'7.31"
$ws.Range("E20").Value = "  -3.25%  "

# Row 21
$ws.Range("D21").Value = "'0.0000106"
$ws.Range("E21").Value = "  -1.99%  "

# Row 22
$ws.Range("D22").Value = "'3.67"
$ws.Range("E22").Value = "  +2.33%  "

# Row 23
$ws.Range("D23").Value = "'73.16"
$ws.Range("E23").Value = "  -2.48%  "

# Row 24
$ws.Range("D24").Value = "'260.82"
$ws.Range("E24").Value = "  -1.57%  "

# Row 25
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = "  -2.90%  "

# Row 26
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$ws.Range("D27").Value = "'11.36"
$ws.Range("E27").Value = "  -0.62%  "

# Row 28
$ws.Range("D28").Value = "'7.41"
$ws.Range("E28").Value = "  -2.77%  "

# Row 29
$ws.Range("E29").Value = "  -2.20%  "

# Row 30
$ws.Range("D30").Value = "'22.36"
$ws.Range("E30").Value = "  -1.88%  "

# Row 31
$ws.Range("D31").Value = "'0.0959"
$ws.Range("E31").Value = "  -3.40%  "

# Row 32
$ws.Range("D32").Value = "'37.01"
$ws.Range("E32").Value = "  -9.38%  "

# Row 33
$ws.Range("D33").Value = "'167.98"
$ws.Range("E33").Value = "  -3.02%  "

# Row 34
$ws.Range("D34").Value = "'2.82"
$ws.Range("E34").Value = "  -4.42%  "

# Row 35
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.132"
$ws.Range("E35").Value = "  -0.64%  "

# Row 36
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'3.28"
$ws.Range("E36").Value = "  +4.87%  "

# Row 37
$ws.Range("D37").Value = "'0.117"
$ws.Range("E37").Value = "  -1.91%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.74"
$ws.Range("E38").Value = "  -5.43%  "

# Row 39
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.97"
$ws.Range("E39").Value = "  +10.64%  "

# Row 40
$ws.Range("D40").Value = "'4.02"
$ws.Range("E40").Value = "  -6.37%  "

# Row 41
$ws.Range("D41").Value = "'0.0354"
$ws.Range("E41").Value = "  -3.53%  "

# Row 42
$ws.Range("D42").Value = "'97.41"
$ws.Range("E42").Value = "  -2.30%  "

# Row 43
$ws.Range("D43").Value = "'70.15"
$ws.Range("E43").Value = "  -2.71%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.885.29"
$ws.Range("E44").Value = "  +14.21%  "

# Row 45
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.230"
$ws.Range("E45").Value = "  -4.67%  "

# Row 46
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").Value = "'6.10"
$ws.Range("E46").Value = "  +4.40%  "

# Row 47
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").Value = "'12.98"
$ws.Range("E47").Value = "  -5.94%  "

# Row 48
$ws.Range("D48").Value = "'86.27"
$ws.Range("E48").Value = "  -0.04%  "

# Row 49
$ws.Range("E49").Value = "  +0.32%  "

# Row 50
$ws.Range("D50").Value = "'112.48"
$ws.Range("E50").Value = "  -3.63%  "

# Row 51
$ws.Range("D51").Value = "'9.31"
$ws.Range("E51").Value = "  -3.00%  "
